$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data matrix for rows 2-10 (A:T), reflecting the addition of the "ECs"
# sending-cluster rows (Dr Hou's advice) ahead of FAPs / sCs.
$data = @(
        @("ECs", "Angpt1", "Tie1", "ECs", 2, 0.6666666666666666, 0.1659303333333333, 0.497791, 0.01127931880021135, 0.01127931880021135, 3, 1, 34.58696433333333, 103.760893, 0.9801058085769325, 0.9801058085769326, 5.73902652081811, 51.651238687363, 0.01105492587287814, 0.01105492587287814),
        @("ECs", "Angpt1", "Tie1", "FAPs", 2, 0.6666666666666666, 0.1659303333333333, 0.497791, 0.01127931880021135, 0.01127931880021135, 3, 1, 0.2577526666666667, 0.773258, 0.007304049101896046, 0.007304049101896047, 0.04276898589755555, 0.384920873078, 0.00008238469835268292, 0.00008238469835268292),
        @("ECs", "Angpt1", "Tie1", "sCs", 2, 0.6666666666666666, 0.1659303333333333, 0.497791, 0.01127931880021135, 0.01127931880021135, 3, 1, 0.4442936666666666, 1.332881, 0.01259014232117133, 0.01259014232117133, 0.07372179620788888, 0.663496165871, 0.0001420082289805244, 0.0001420082289805244),
        @("FAPs", "Angpt1", "Tie1", "ECs", 3, 1, 12.95459633333333, 38.863789, 0.8806046431437031, 0.8806046431437031, 3, 1, 34.58696433333333, 103.760893, 0.9801058085769325, 0.9801058085769326, 448.0601613337308, 4032.541452003577, 0.8630857258049602, 0.8630857258049603),
        @("FAPs", "Angpt1", "Tie1", "FAPs", 3, 1, 12.95459633333333, 38.863789, 0.8806046431437031, 0.8806046431437031, 3, 1, 0.2577526666666667, 0.773258, 0.007304049101896046, 0.007304049101896047, 3.339081750506889, 30.051735754562, 0.006431979552879253, 0.006431979552879254),
        @("FAPs", "Angpt1", "Tie1", "sCs", 3, 1, 12.95459633333333, 38.863789, 0.8806046431437031, 0.8806046431437031, 3, 1, 0.4442936666666666, 1.332881, 0.01259014232117133, 0.01259014232117133, 5.755645105123222, 51.80080594610899, 0.01108693778586352, 0.01108693778586352),
        @("sCs", "Angpt1", "Tie1", "ECs", 3, 1, 1.590497666666667, 4.771493, 0.1081160380560855, 0.1081160380560855, 3, 1, 34.58696433333333, 103.760893, 0.9801058085769325, 0.9801058085769326, 55.01048606924989, 495.0943746232489, 0.1059651568990941, 0.1059651568990941),
        @("sCs", "Angpt1", "Tie1", "FAPs", 3, 1, 1.590497666666667, 4.771493, 0.1081160380560855, 0.1081160380560855, 3, 1, 0.2577526666666667, 0.773258, 0.007304049101896046, 0.007304049101896047, 0.4099550149104444, 3.689595134194, 0.0007896848506641101, 0.0007896848506641102),
        @("sCs", "Angpt1", "Tie1", "sCs", 3, 1, 1.590497666666667, 4.771493, 0.1081160380560855, 0.1081160380560855, 3, 1, 0.4442936666666666, 1.332881, 0.01259014232117133, 0.01259014232117133, 0.706648040148111, 6.359832361332999, 0.001361196306327293, 0.001361196306327293)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
